# Replace each arithmetic-equation table cell's text with its new value.
# Every "old" equation string below is unique within the document, so a
# simple literal Find/Replace (wildcards off) for each pair is safe and
# unambiguous.
$d = $word.ActiveDocument

$d.Content.Find.Execute("76+9=85", $false, $false, $false, $false, $false, $true, 1, $false, "40+2=42", 2) | Out-Null
$d.Content.Find.Execute("56-7=49", $false, $false, $false, $false, $false, $true, 1, $false, "51-33=18", 2) | Out-Null
$d.Content.Find.Execute("59-27=32", $false, $false, $false, $false, $false, $true, 1, $false, "38-33=5", 2) | Out-Null
$d.Content.Find.Execute("55-35=20", $false, $false, $false, $false, $false, $true, 1, $false, "64+18=82", 2) | Out-Null
$d.Content.Find.Execute("20+5=25", $false, $false, $false, $false, $false, $true, 1, $false, "46-36=10", 2) | Out-Null
$d.Content.Find.Execute("14+60=74", $false, $false, $false, $false, $false, $true, 1, $false, "84-43=41", 2) | Out-Null
$d.Content.Find.Execute("37+10=47", $false, $false, $false, $false, $false, $true, 1, $false, "40-14=26", 2) | Out-Null
$d.Content.Find.Execute("33+28=61", $false, $false, $false, $false, $false, $true, 1, $false, "87-46=41", 2) | Out-Null
$d.Content.Find.Execute("80-34=46", $false, $false, $false, $false, $false, $true, 1, $false, "26-17=9", 2) | Out-Null
$d.Content.Find.Execute("1+53=54", $false, $false, $false, $false, $false, $true, 1, $false, "84+13=97", 2) | Out-Null
$d.Content.Find.Execute("38+28=66", $false, $false, $false, $false, $false, $true, 1, $false, "8+46=54", 2) | Out-Null
$d.Content.Find.Execute("12-8=4", $false, $false, $false, $false, $false, $true, 1, $false, "88+1=89", 2) | Out-Null
$d.Content.Find.Execute("2+67=69", $false, $false, $false, $false, $false, $true, 1, $false, "57+25=82", 2) | Out-Null
$d.Content.Find.Execute("45+7=52", $false, $false, $false, $false, $false, $true, 1, $false, "41-14=27", 2) | Out-Null
$d.Content.Find.Execute("69+24=93", $false, $false, $false, $false, $false, $true, 1, $false, "96-33=63", 2) | Out-Null
$d.Content.Find.Execute("82-6=76", $false, $false, $false, $false, $false, $true, 1, $false, "52-0=52", 2) | Out-Null
$d.Content.Find.Execute("8+58=66", $false, $false, $false, $false, $false, $true, 1, $false, "13+0=13", 2) | Out-Null
$d.Content.Find.Execute("85-39=46", $false, $false, $false, $false, $false, $true, 1, $false, "10+73=83", 2) | Out-Null
$d.Content.Find.Execute("62-17=45", $false, $false, $false, $false, $false, $true, 1, $false, "37+54=91", 2) | Out-Null
$d.Content.Find.Execute("76+23=99", $false, $false, $false, $false, $false, $true, 1, $false, "28+51=79", 2) | Out-Null
$d.Content.Find.Execute("29+40=69", $false, $false, $false, $false, $false, $true, 1, $false, "24+52=76", 2) | Out-Null
$d.Content.Find.Execute("97-4=93", $false, $false, $false, $false, $false, $true, 1, $false, "11+3=14", 2) | Out-Null
$d.Content.Find.Execute("66+15=81", $false, $false, $false, $false, $false, $true, 1, $false, "93-75=18", 2) | Out-Null
$d.Content.Find.Execute("55+35=90", $false, $false, $false, $false, $false, $true, 1, $false, "54+30=84", 2) | Out-Null
$d.Content.Find.Execute("53+44=97", $false, $false, $false, $false, $false, $true, 1, $false, "49-41=8", 2) | Out-Null
$d.Content.Find.Execute("89-88=1", $false, $false, $false, $false, $false, $true, 1, $false, "57+28=85", 2) | Out-Null
$d.Content.Find.Execute("66-57=9", $false, $false, $false, $false, $false, $true, 1, $false, "12+50=62", 2) | Out-Null
$d.Content.Find.Execute("94-75=19", $false, $false, $false, $false, $false, $true, 1, $false, "46+51=97", 2) | Out-Null
$d.Content.Find.Execute("98-95=3", $false, $false, $false, $false, $false, $true, 1, $false, "81-9=72", 2) | Out-Null
$d.Content.Find.Execute("28+26=54", $false, $false, $false, $false, $false, $true, 1, $false, "15+32=47", 2) | Out-Null
$d.Content.Find.Execute("26-9=17", $false, $false, $false, $false, $false, $true, 1, $false, "96-28=68", 2) | Out-Null
$d.Content.Find.Execute("87-74=13", $false, $false, $false, $false, $false, $true, 1, $false, "36+10=46", 2) | Out-Null
$d.Content.Find.Execute("47-47=0", $false, $false, $false, $false, $false, $true, 1, $false, "8+83=91", 2) | Out-Null
$d.Content.Find.Execute("30+23=53", $false, $false, $false, $false, $false, $true, 1, $false, "81-13=68", 2) | Out-Null
$d.Content.Find.Execute("53+36=89", $false, $false, $false, $false, $false, $true, 1, $false, "89+8=97", 2) | Out-Null
$d.Content.Find.Execute("14+28=42", $false, $false, $false, $false, $false, $true, 1, $false, "93-69=24", 2) | Out-Null
$d.Content.Find.Execute("85-84=1", $false, $false, $false, $false, $false, $true, 1, $false, "58+26=84", 2) | Out-Null
$d.Content.Find.Execute("65+5=70", $false, $false, $false, $false, $false, $true, 1, $false, "32+9=41", 2) | Out-Null
$d.Content.Find.Execute("1+9=10", $false, $false, $false, $false, $false, $true, 1, $false, "24+46=70", 2) | Out-Null
$d.Content.Find.Execute("60-1=59", $false, $false, $false, $false, $false, $true, 1, $false, "3+53=56", 2) | Out-Null
$d.Content.Find.Execute("17+74=91", $false, $false, $false, $false, $false, $true, 1, $false, "20+27=47", 2) | Out-Null
$d.Content.Find.Execute("2+38=40", $false, $false, $false, $false, $false, $true, 1, $false, "11+72=83", 2) | Out-Null
$d.Content.Find.Execute("67-23=44", $false, $false, $false, $false, $false, $true, 1, $false, "73+11=84", 2) | Out-Null
$d.Content.Find.Execute("70+1=71", $false, $false, $false, $false, $false, $true, 1, $false, "14+49=63", 2) | Out-Null
$d.Content.Find.Execute("70-58=12", $false, $false, $false, $false, $false, $true, 1, $false, "84-20=64", 2) | Out-Null
$d.Content.Find.Execute("57-36=21", $false, $false, $false, $false, $false, $true, 1, $false, "29+39=68", 2) | Out-Null
$d.Content.Find.Execute("35+23=58", $false, $false, $false, $false, $false, $true, 1, $false, "36+11=47", 2) | Out-Null
$d.Content.Find.Execute("80-12=68", $false, $false, $false, $false, $false, $true, 1, $false, "84-62=22", 2) | Out-Null
$d.Content.Find.Execute("6+49=55", $false, $false, $false, $false, $false, $true, 1, $false, "6+28=34", 2) | Out-Null
$d.Content.Find.Execute("79-19=60", $false, $false, $false, $false, $false, $true, 1, $false, "93-84=9", 2) | Out-Null
$d.Content.Find.Execute("84-53=31", $false, $false, $false, $false, $false, $true, 1, $false, "29+56=85", 2) | Out-Null
$d.Content.Find.Execute("62+12=74", $false, $false, $false, $false, $false, $true, 1, $false, "96-87=9", 2) | Out-Null
$d.Content.Find.Execute("76-70=6", $false, $false, $false, $false, $false, $true, 1, $false, "35-6=29", 2) | Out-Null
$d.Content.Find.Execute("82-16=66", $false, $false, $false, $false, $false, $true, 1, $false, "38+40=78", 2) | Out-Null
$d.Content.Find.Execute("79-67=12", $false, $false, $false, $false, $false, $true, 1, $false, "43-1=42", 2) | Out-Null
$d.Content.Find.Execute("28+45=73", $false, $false, $false, $false, $false, $true, 1, $false, "13+65=78", 2) | Out-Null
$d.Content.Find.Execute("70-19=51", $false, $false, $false, $false, $false, $true, 1, $false, "37+61=98", 2) | Out-Null
$d.Content.Find.Execute("0+83=83", $false, $false, $false, $false, $false, $true, 1, $false, "43+30=73", 2) | Out-Null
$d.Content.Find.Execute("3+76=79", $false, $false, $false, $false, $false, $true, 1, $false, "81+6=87", 2) | Out-Null
$d.Content.Find.Execute("30+54=84", $false, $false, $false, $false, $false, $true, 1, $false, "91-39=52", 2) | Out-Null
$d.Content.Find.Execute("59-13=46", $false, $false, $false, $false, $false, $true, 1, $false, "10+14=24", 2) | Out-Null
$d.Content.Find.Execute("59-9=50", $false, $false, $false, $false, $false, $true, 1, $false, "91-34=57", 2) | Out-Null
$d.Content.Find.Execute("34+5=39", $false, $false, $false, $false, $false, $true, 1, $false, "71-29=42", 2) | Out-Null
$d.Content.Find.Execute("12-6=6", $false, $false, $false, $false, $false, $true, 1, $false, "82-62=20", 2) | Out-Null
$d.Content.Find.Execute("74-30=44", $false, $false, $false, $false, $false, $true, 1, $false, "9+79=88", 2) | Out-Null
$d.Content.Find.Execute("52+32=84", $false, $false, $false, $false, $false, $true, 1, $false, "99-36=63", 2) | Out-Null
$d.Content.Find.Execute("75-63=12", $false, $false, $false, $false, $false, $true, 1, $false, "64+13=77", 2) | Out-Null
$d.Content.Find.Execute("33+6=39", $false, $false, $false, $false, $false, $true, 1, $false, "51-36=15", 2) | Out-Null
$d.Content.Find.Execute("41-17=24", $false, $false, $false, $false, $false, $true, 1, $false, "86-9=77", 2) | Out-Null
$d.Content.Find.Execute("19-11=8", $false, $false, $false, $false, $false, $true, 1, $false, "67-30=37", 2) | Out-Null
$d.Content.Find.Execute("9+90=99", $false, $false, $false, $false, $false, $true, 1, $false, "25-14=11", 2) | Out-Null
$d.Content.Find.Execute("70-54=16", $false, $false, $false, $false, $false, $true, 1, $false, "46-26=20", 2) | Out-Null
$d.Content.Find.Execute("92-43=49", $false, $false, $false, $false, $false, $true, 1, $false, "27+17=44", 2) | Out-Null
$d.Content.Find.Execute("9-7=2", $false, $false, $false, $false, $false, $true, 1, $false, "52-38=14", 2) | Out-Null
$d.Content.Find.Execute("87+12=99", $false, $false, $false, $false, $false, $true, 1, $false, "94-93=1", 2) | Out-Null
$d.Content.Find.Execute("83-9=74", $false, $false, $false, $false, $false, $true, 1, $false, "63+35=98", 2) | Out-Null
$d.Content.Find.Execute("84-55=29", $false, $false, $false, $false, $false, $true, 1, $false, "52-50=2", 2) | Out-Null
$d.Content.Find.Execute("96-94=2", $false, $false, $false, $false, $false, $true, 1, $false, "15+32=47", 2) | Out-Null
$d.Content.Find.Execute("14+72=86", $false, $false, $false, $false, $false, $true, 1, $false, "5+44=49", 2) | Out-Null
$d.Content.Find.Execute("16+78=94", $false, $false, $false, $false, $false, $true, 1, $false, "47+43=90", 2) | Out-Null
$d.Content.Find.Execute("65-40=25", $false, $false, $false, $false, $false, $true, 1, $false, "34+43=77", 2) | Out-Null
$d.Content.Find.Execute("87-64=23", $false, $false, $false, $false, $false, $true, 1, $false, "84-43=41", 2) | Out-Null
$d.Content.Find.Execute("33+35=68", $false, $false, $false, $false, $false, $true, 1, $false, "54+44=98", 2) | Out-Null
$d.Content.Find.Execute("21-15=6", $false, $false, $false, $false, $false, $true, 1, $false, "26+13=39", 2) | Out-Null
$d.Content.Find.Execute("4+84=88", $false, $false, $false, $false, $false, $true, 1, $false, "86-65=21", 2) | Out-Null
$d.Content.Find.Execute("69+19=88", $false, $false, $false, $false, $false, $true, 1, $false, "23-5=18", 2) | Out-Null
$d.Content.Find.Execute("85-0=85", $false, $false, $false, $false, $false, $true, 1, $false, "32+6=38", 2) | Out-Null
$d.Content.Find.Execute("39-18=21", $false, $false, $false, $false, $false, $true, 1, $false, "66-36=30", 2) | Out-Null
$d.Content.Find.Execute("98-43=55", $false, $false, $false, $false, $false, $true, 1, $false, "8+22=30", 2) | Out-Null
$d.Content.Find.Execute("8+52=60", $false, $false, $false, $false, $false, $true, 1, $false, "47+18=65", 2) | Out-Null
$d.Content.Find.Execute("57+39=96", $false, $false, $false, $false, $false, $true, 1, $false, "51+38=89", 2) | Out-Null
$d.Content.Find.Execute("0+72=72", $false, $false, $false, $false, $false, $true, 1, $false, "72-71=1", 2) | Out-Null
$d.Content.Find.Execute("29+34=63", $false, $false, $false, $false, $false, $true, 1, $false, "32-9=23", 2) | Out-Null
$d.Content.Find.Execute("60+39=99", $false, $false, $false, $false, $false, $true, 1, $false, "28+27=55", 2) | Out-Null
$d.Content.Find.Execute("5+89=94", $false, $false, $false, $false, $false, $true, 1, $false, "54+39=93", 2) | Out-Null
$d.Content.Find.Execute("7+54=61", $false, $false, $false, $false, $false, $true, 1, $false, "57-21=36", 2) | Out-Null
$d.Content.Find.Execute("9+3=12", $false, $false, $false, $false, $false, $true, 1, $false, "32-24=8", 2) | Out-Null
$d.Content.Find.Execute("69+25=94", $false, $false, $false, $false, $false, $true, 1, $false, "53+46=99", 2) | Out-Null
$d.Content.Find.Execute("85-56=29", $false, $false, $false, $false, $false, $true, 1, $false, "51-0=51", 2) | Out-Null
$d.Content.Find.Execute("58-40=18", $false, $false, $false, $false, $false, $true, 1, $false, "88+9=97", 2) | Out-Null
